$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) edits to be stored as text, matching the
# original inline-string cell type (avoids Excel auto-parsing values
# like "298.76" or "1.001" as numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.458.56"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.646.18"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "298.76"
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3788"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3540"
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.03"
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08088"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.216"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.389"
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.332"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001198"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.651.60"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.07"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06945"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.767"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.37"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.43"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.455.54"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.494"
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.900"
$ws.Range("E26").Value = "  -6.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.86"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.15"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.205"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.72"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.831.89"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.925"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.147"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.49"
$ws.Range("E34").Value = "  -3.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9903"
$ws.Range("E35").Value = "  -8.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02712"
$ws.Range("E36").Value = "  -4.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08748"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2435"
$ws.Range("E38").Value = "  -3.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.920"
$ws.Range("E39").Value = "  -3.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "13.00"
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06786"
$ws.Range("E41").Value = "  -3.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6884"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.294"
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.64"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6352"
$ws.Range("E46").Value = "  -3.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.250"
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.908"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07721"
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.42"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.148"
$ws.Range("E51").Value = "  -4.22%  "
